# Adds the 27/10/2024 diary entry (friend-request rejection note, roadmap
# bullet list, and closing remarks) to the end of the document body, right
# after the 25/10/2024 paragraph and before the sectPr.

$d = $word.ActiveDocument

function New-PkgXml([string]$innerBodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerBodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

function Escape-Xml([string]$text) {
    return $text -replace '&','&amp;' -replace '<','&lt;' -replace '>','&gt;' -replace '"','&quot;'
}

# Plain "diary" paragraph: ind left=708, purple (6600FF) text, no list numbering.
function Add-PlainPara($cursorRange, [string]$text) {
    $t = Escape-Xml $text
    $p = '<w:p><w:pPr><w:ind w:left="708"/><w:rPr><w:color w:val="6600FF"/></w:rPr></w:pPr>' +
         '<w:r><w:rPr><w:color w:val="6600FF"/></w:rPr><w:t>' + $t + '</w:t></w:r></w:p>'
    $cursorRange.InsertXML((New-PkgXml $p))
}

# Roadmap bullet paragraph: "Prrafodelista" style, ilvl=3 / numId=3 list, purple text.
function Add-ListPara($cursorRange, [string]$text) {
    $t = Escape-Xml $text
    $p = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="3"/></w:numPr>' +
         '<w:rPr><w:color w:val="6600FF"/></w:rPr></w:pPr>' +
         '<w:r><w:rPr><w:color w:val="6600FF"/></w:rPr><w:t>' + $t + '</w:t></w:r></w:p>'
    $cursorRange.InsertXML((New-PkgXml $p))
}

# Start a fresh trailing paragraph (inherits formatting of the current last
# paragraph) that will act as our insertion cursor.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$cursor = $d.Paragraphs.Last.Range

Add-PlainPara $cursor "27/10/2024: Añado el método para rechazar peticiones de amistad. Solo faltan comentarios para algunas clases, los iré añadiendo según tenga tiempo. Cosas a implementar que creo necesarias e iré haciendo poco a poco."
$cursor = $d.Paragraphs.Last.Range

Add-ListPara $cursor "Primero: que los usuarios tengan una foto de perfil y puedan subir una"
$cursor = $d.Paragraphs.Last.Range

Add-ListPara $cursor "Segundo: una activity que muestre los datos de cada usuario y los eventos a los que está unido y coincidan con los que tiene el usuario"
$cursor = $d.Paragraphs.Last.Range

Add-ListPara $cursor "Tercero: una forma de eliminar amigos"
$cursor = $d.Paragraphs.Last.Range

Add-ListPara $cursor "Cuarto: añadir al login un “recuérdame” para que un usuario no tenga que logearse varias veces."
$cursor = $d.Paragraphs.Last.Range

Add-ListPara $cursor "Quinto: poner al usuario en la appbar y que pueda hacer logout"
$cursor = $d.Paragraphs.Last.Range

Add-ListPara $cursor "Sexto: añadir la appbar a todas las activities"
$cursor = $d.Paragraphs.Last.Range

Add-ListPara $cursor "Séptimo: unificar los estilos de todas las activities"
$cursor = $d.Paragraphs.Last.Range

Add-PlainPara $cursor "No hay una preferencia como tal en este orden, pero son las cosas que me faltarían de momento por implementar. Creo que deberíamos centrar esfuerzos en los eventos para ya poder ir concretando al app."
$cursor = $d.Paragraphs.Last.Range

Add-PlainPara $cursor "--IMPORTANTE—cuando tenga los comentarios de las clases subiré todo lo que tengo hasta ahora al main para que tengamos todas las clases como la de usuario disponibles todos."

# The final empty paragraph (ind left=708, purple rPr, no run) is whatever is
# left over as the new trailing cursor paragraph - leave it untouched.

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
